# Cypress Ascendant Services LLC sample template - bulk excel template fixes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the sample data row: "Full Searchs" -> "Full Search"
$ws.Range("F2").Value = "Full Search"

# Leave the sheet's active-cell selection on F2 (where the fix was made)
$ws.Range("F2").Select()
